# Fruta / hortaliza, semanal
# Insert two new weekly price rows (213 and 214) for "Vega Modelo de Temuco - Durazno",
# pushing the existing rows 213-223 down to 215-225.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 213, shifting rows 213:223 down to 215:225.
$ws.Rows("213:214").Insert()

# New row 213 - Kurakata, Primera
$ws.Range("A213").Value = 10
$ws.Range("B213").Value = "Vega Modelo de Temuco"
$ws.Range("C213").Value = "La Araucanía"
$ws.Range("D213").Value = 44585
$ws.Range("E213").Value = 9
$ws.Range("F213").Value = "Fruta"
$ws.Range("G213").Value = 100103
$ws.Range("H213").Value = "Frutos de hueso (carozo)"
$ws.Range("I213").Value = 100103004
$ws.Range("J213").Value = "Durazno"
$ws.Range("K213").Value = "Kurakata"
$ws.Range("L213").Value = "Primera"
$ws.Range("M213").Value = 130
$ws.Range("N213").Value = 16000
$ws.Range("O213").Value = 17000
$ws.Range("P213").Value = 16615
$ws.Range("Q213").Value = "$/bandeja 18 kilos granel"
$ws.Range("R213").Value = "Región de O'Higgins"
$ws.Range("S213").Value = 923
$ws.Range("T213").Value = 18

# New row 214 - Kurakata, Primera (bins)
$ws.Range("A214").Value = 10
$ws.Range("B214").Value = "Vega Modelo de Temuco"
$ws.Range("C214").Value = "La Araucanía"
$ws.Range("D214").Value = 44585
$ws.Range("E214").Value = 9
$ws.Range("F214").Value = "Fruta"
$ws.Range("G214").Value = 100103
$ws.Range("H214").Value = "Frutos de hueso (carozo)"
$ws.Range("I214").Value = 100103004
$ws.Range("J214").Value = "Durazno"
$ws.Range("K214").Value = "Kurakata"
$ws.Range("L214").Value = "Primera"
$ws.Range("M214").Value = 10
$ws.Range("N214").Value = 350000
$ws.Range("O214").Value = 350000
$ws.Range("P214").Value = 350000
$ws.Range("Q214").Value = "$/bins (400 kilos)"
$ws.Range("R214").Value = "Región de O'Higgins"
$ws.Range("S214").Value = 875
$ws.Range("T214").Value = 400

Write-Host "Inserted new rows 213-214 for Vega Modelo de Temuco - Durazno"
